$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at row 3 (EmailNotification) - existing rows 3,4,5 shift
#    down to become rows 4,5,6.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "EmailNotification"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "Suited to Manual"
$ws.Range("E3").Value = "Email"

# Copy the cell formatting (styles) from row 2 so the new row looks like the
# rest of the table (same borders / alignment / number formats per column).
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B2:D2").Copy()
$ws.Range("B3:D3").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the conditional formatting so the newly inserted row gets its own
#    set of rules (mirroring the ones already present on the D column) while
#    re-numbering the priorities of the pre-existing rules to make room.
# ---------------------------------------------------------------------------
$existing = $ws.Range("D1:D1048576").FormatConditions
$existingItems = @()
for ($i = 1; $i -le $existing.Count; $i++) {
    $existingItems += $existing.Item($i)
}
$newPriorities = @(17, 9, 10, 11, 12, 13, 14, 15)
for ($i = 0; $i -lt $existingItems.Length; $i++) {
    $existingItems[$i].Priority = $newPriorities[$i]
}

$d3 = $ws.Range("D3").FormatConditions

$blank = $d3.Add(2, 0, "LEN(TRIM(D3))>0")
$blank.StopIfTrue = $true
$blank.Priority = 8
$blank.Interior.Color = 255

$finished = $d3.Add(9, 0, "Finished", $null, "Finished")
$finished.StopIfTrue = $true
$finished.Priority = 1
$finished.Interior.Color = 6750054

$automated = $d3.Add(9, 0, "Automated", $null, "Automated")
$automated.StopIfTrue = $true
$automated.Priority = 2
$automated.Interior.Color = 6750054

$underReview = $d3.Add(9, 0, "Under Review", $null, "Under Review")
$underReview.StopIfTrue = $true
$underReview.Priority = 3
$underReview.Interior.Color = 5287936

$testing = $d3.Add(9, 0, "Testing", $null, "Testing")
$testing.StopIfTrue = $true
$testing.Priority = 4
$testing.Interior.Color = 5287936

$writing = $d3.Add(9, 0, "Writing", $null, "Writing")
$writing.StopIfTrue = $true
$writing.Priority = 5
$writing.Interior.Color = 15773696

$readyToWrite = $d3.Add(9, 0, "Ready to Write", $null, "Ready to Write")
$readyToWrite.StopIfTrue = $true
$readyToWrite.Priority = 6
$readyToWrite.Interior.Color = 15773696

$unwrittenKeywords = $d3.Add(9, 0, "Unwritten Keywords", $null, "Unwritten Keywords")
$unwrittenKeywords.StopIfTrue = $true
$unwrittenKeywords.Priority = 7
$unwrittenKeywords.Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. Selection moved to E4, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("E4").Select()

$excel.Calculate()
